$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 4497.5
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 5495
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 5495
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -6743

# ALC row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 4497.5
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 5495
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 27475
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -33715

# ALC row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 5098.5713
$ws.Range("I86").Value = 2404.2856
$ws.Range("J86").Value = 7792.857
$ws.Range("K86").Value = 2404.2856
$ws.Range("L86").Value = 7792.857
$ws.Range("M86").Value = -1281.2856
$ws.Range("N86").Value = -10038.857

# ALC row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 5098.5713
$ws.Range("I89").Value = 2404.2856
$ws.Range("J89").Value = 7792.857
$ws.Range("K89").Value = 12021.428
$ws.Range("L89").Value = 38964.285
$ws.Range("M89").Value = -6405.428
$ws.Range("N89").Value = -50196.285

# ALC row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 829.4138
$ws.Range("I107").Value = 829.9
$ws.Range("J107").Value = 828.3333
$ws.Range("K107").Value = 829.9
$ws.Range("L107").Value = 828.3333
$ws.Range("M107").Value = 1090.1
$ws.Range("N107").Value = -4668.3333

# ALC row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 821.4888999999999
$ws.Range("J129").Value = 849.2683
$ws.Range("L129").Value = 2547.8049
$ws.Range("N129").Value = -12547.8049

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 21740564
$ws.Range("I137").Value = 977.9697
$ws.Range("J137").Value = 76925660
$ws.Range("K137").Value = 2933.9091
$ws.Range("L137").Value = 230776980
$ws.Range("M137").Value = -383.9090999999999
$ws.Range("N137").Value = -230782080

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2723.878
$ws.Range("I138").Value = 2029.75
$ws.Range("J138").Value = 3168.12
$ws.Range("K138").Value = 6089.25
$ws.Range("L138").Value = 9504.360000000001
$ws.Range("M138").Value = -949.25
$ws.Range("N138").Value = -19784.36

$ws = $wb.Worksheets.Item("ARM")
# ARM row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 29212.5
$ws.Range("J139").Value = 29212.5
$ws.Range("L139").Value = 29212.5
$ws.Range("N139").Value = -39492.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1177.4445
$ws.Range("I86").Value = 1090.2858
$ws.Range("J86").Value = 1482.5
$ws.Range("K86").Value = 1090.2858
$ws.Range("L86").Value = 1482.5
$ws.Range("M86").Value = 32.71419999999989
$ws.Range("N86").Value = -3728.5

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1177.4445
$ws.Range("I89").Value = 1090.2858
$ws.Range("J89").Value = 1482.5
$ws.Range("K89").Value = 5451.429
$ws.Range("L89").Value = 7412.5
$ws.Range("M89").Value = 164.5709999999999
$ws.Range("N89").Value = -18644.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2916.8906
$ws.Range("I58").Value = 1217.909
$ws.Range("K58").Value = 1217.909
$ws.Range("M58").Value = -1014.909

# CRP row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 13753
$ws.Range("I122").Value = 17520.666
$ws.Range("K122").Value = 52561.99800000001
$ws.Range("M122").Value = -50111.99800000001

# CRP row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2640.0833
$ws.Range("I132").Value = 2390.1853
$ws.Range("J132").Value = 3389.7778
$ws.Range("K132").Value = 7170.5559
$ws.Range("L132").Value = 10169.3334
$ws.Range("M132").Value = -4640.5559
$ws.Range("N132").Value = -15229.3334

# CRP row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2916.8906
$ws.Range("I136").Value = 1217.909
$ws.Range("K136").Value = 3653.727
$ws.Range("M136").Value = -1103.727

$ws = $wb.Worksheets.Item("CUL")
# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 5475.84
$ws.Range("J131").Value = 3908.7368
$ws.Range("L131").Value = 11726.2104
$ws.Range("N131").Value = -21806.2104

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 11178.546
$ws.Range("I80").Value = 2457.8572
$ws.Range("J80").Value = 26439.75
$ws.Range("K80").Value = 2457.8572
$ws.Range("L80").Value = 26439.75
$ws.Range("M80").Value = -1459.8572
$ws.Range("N80").Value = -28435.75

# GSM row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 11178.546
$ws.Range("I83").Value = 2457.8572
$ws.Range("J83").Value = 26439.75
$ws.Range("K83").Value = 12289.286
$ws.Range("L83").Value = 132198.75
$ws.Range("M83").Value = -7297.286
$ws.Range("N83").Value = -142182.75

# GSM row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 2157.0588
$ws.Range("I97").Value = 2179.375
$ws.Range("J97").Value = 1800
$ws.Range("K97").Value = 2179.375
$ws.Range("L97").Value = 1800
$ws.Range("M97").Value = -1683.375
$ws.Range("N97").Value = -2792

# GSM row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1981.8182
$ws.Range("I122").Value = 2237.5
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 6712.5
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -4262.5
$ws.Range("N122").Value = -8800

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 2814.0833
$ws.Range("I40").Value = 3173.5
$ws.Range("J40").Value = 2634.375
$ws.Range("K40").Value = 3173.5
$ws.Range("L40").Value = 2634.375
$ws.Range("M40").Value = -3037.5
$ws.Range("N40").Value = -2906.375

# LTW row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 1584.2106
$ws.Range("I82").Value = 1514.2858
$ws.Range("J82").Value = 1780
$ws.Range("K82").Value = 1514.2858
$ws.Range("L82").Value = 1780
$ws.Range("M82").Value = -1153.2858
$ws.Range("N82").Value = -2502

# LTW row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 1584.2106
$ws.Range("I85").Value = 1514.2858
$ws.Range("J85").Value = 1780
$ws.Range("K85").Value = 1514.2858
$ws.Range("L85").Value = 1780
$ws.Range("M85").Value = -266.2858000000001
$ws.Range("N85").Value = -4276

# LTW row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 2681.4736
$ws.Range("I100").Value = 1855.7858
$ws.Range("J100").Value = 4993.4
$ws.Range("K100").Value = 1855.7858
$ws.Range("L100").Value = 4993.4
$ws.Range("M100").Value = -1314.7858
$ws.Range("N100").Value = -6075.4

$ws = $wb.Worksheets.Item("WVR")
# WVR row 96: Skills on Display / Ruby Cotton Cloth
$ws.Range("H96").Value = 9107687
$ws.Range("I96").Value = 25003476
$ws.Range("J96").Value = 24380.285
$ws.Range("K96").Value = 25003476
$ws.Range("L96").Value = 24380.285
$ws.Range("M96").Value = -25002103
$ws.Range("N96").Value = -27126.285

# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 4258.104
$ws.Range("I132").Value = 4829.1055
$ws.Range("J132").Value = 2088.3
$ws.Range("K132").Value = 14487.3165
$ws.Range("L132").Value = 6264.900000000001
$ws.Range("M132").Value = -11957.3165
$ws.Range("N132").Value = -11324.9

# WVR row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 8114.34
$ws.Range("I136").Value = 10542.757
$ws.Range("J136").Value = 1202.6923
$ws.Range("K136").Value = 31628.271
$ws.Range("L136").Value = 3608.0769
$ws.Range("M136").Value = -29078.271
$ws.Range("N136").Value = -8708.0769
